{"js": "// Remove the stray \"_GoBack\" bookmark (it will be re-added below, after the\n// new acceptance-criteria text) and insert the two acceptance-criteria\n// lines into the empty paragraph that follows \"ACCEPTANCE CRITERIA\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1. Delete the existing \"_GoBack\" bookmark that currently sits in the\n//    USER STORY paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Locate the empty paragraph right after the \"ACCEPTANCE CRITERIA\"\n//    heading paragraph.\nconst items = paragraphs.items;\nlet targetParagraph = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"ACCEPTANCE CRITERIA\") !== -1) {\n    targetParagraph = items[i + 1];\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not locate paragraph following 'ACCEPTANCE CRITERIA'.\");\n}\n\n// 3. Insert the two acceptance-criteria sentences (black font colour),\n//    separated by a manual line break, as two runs inside that paragraph.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r>' +\n  '<w:rPr><w:color w:val=\"000000\"/></w:rPr>' +\n  '<w:t>1. On entering search parameter and clicking on search button user navigates to List of provider page</w:t>' +\n  '</w:r>' +\n  '<w:r>' +\n  '<w:rPr><w:color w:val=\"000000\"/></w:rPr>' +\n  '<w:br/>' +\n  '<w:t>2. User can refine the search by using different filter criteria Country, city, provider and quality star rating</w:t>' +\n  '</w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntargetParagraph.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n\n// 4. Re-add the \"_GoBack\" bookmark at the end of that paragraph (now after\n//    the newly-inserted text), matching the edited document.\nconst endRange = targetParagraph.getRange(Word.RangeLocation.end);\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Remove the stray \"_GoBack\" bookmark that currently sits in the USER STORY\n# paragraph, and insert the two acceptance-criteria lines (black font\n# colour, separated by a manual line break) into the empty paragraph that\n# follows the \"ACCEPTANCE CRITERIA\" heading, re-adding the \"_GoBack\"\n# bookmark right after the newly inserted text.\n\n$d = $word.ActiveDocument\n\n# 1. Delete the existing \"_GoBack\" bookmark (it currently lives at the end\n#    of the USER STORY paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Locate the empty paragraph right after the \"ACCEPTANCE CRITERIA\"\n#    heading paragraph.\n$target = $null\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*ACCEPTANCE CRITERIA*\") {\n        $targetIndex = $i + 1\n        $target = $d.Paragraphs.Item($targetIndex)\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate paragraph following 'ACCEPTANCE CRITERIA'.\"\n}\n\n# 3. Insert the two acceptance-criteria sentences (black font colour,\n#    joined by a manual line break) as two runs inside that paragraph. A\n#    one-character placeholder run (\"Z\") is appended after them; it is used\n#    below to safely re-create the \"_GoBack\" bookmark and is then removed.\n$tr = $target.Range\n$insertPos = $tr.End - 1\n$insertRange = $d.Range($insertPos, $insertPos)\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r>' +\n    '<w:rPr><w:color w:val=\"000000\"/></w:rPr>' +\n    '<w:t>1. On entering search parameter and clicking on search button user navigates to List of provider page</w:t>' +\n    '</w:r>' +\n    '<w:r>' +\n    '<w:rPr><w:color w:val=\"000000\"/></w:rPr>' +\n    '<w:br/>' +\n    '<w:t>2. User can refine the search by using different filter criteria Country, city, provider and quality star rating</w:t>' +\n    '</w:r>' +\n    '<w:r><w:t>Z</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$insertRange.InsertXML($ooxml) | Out-Null\n\n# 4. Re-add the \"_GoBack\" bookmark. A *collapsed* (zero-length) Range\n#    placed right at a paragraph boundary is unreliable for\n#    Bookmarks.Add in this host, so the bookmark is first created around\n#    the one-character \"Z\" placeholder (a non-collapsed Range, which is\n#    reliable) and the placeholder is deleted immediately afterwards,\n#    leaving the bookmark collapsed at the correct position.\n$target2 = $d.Paragraphs.Item($targetIndex)\n$tr2 = $target2.Range\n$placeholderEnd = $tr2.End - 1\n$placeholderStart = $placeholderEnd - 1\n$bmRange = $d.Range($placeholderStart, $placeholderEnd)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n$delRange = $d.Range($placeholderStart, $placeholderEnd)\n$delRange.Delete()\n"}
